$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in this week's (21/06/2018) actuals in column I ---
$ws.Range("I2").Value = 7
$ws.Range("I4").Value = 15
$ws.Range("I5").Value = 30
$ws.Range("I8").Value = 30
$ws.Range("I10").Value = 13
$ws.Range("I13").Value = 2
$ws.Range("I16").Value = 23
$ws.Range("I17").Value = 15
$ws.Range("I18").Value = 23

# --- 2. Correct two previously mis-entered values ---
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 2

# --- 3. Drop the three future weeks (28/06, 05/07, 12/07) that no longer apply ---
$ws.Columns("J:L").Delete() | Out-Null

# --- 4. Add the weekly-checksum column back in column J ---
$ws.Range("J2").Formula = "=B2-SUM(C2:I2)"
$ws.Range("J3").Formula = "=B3-SUM(C3:I3)"
$ws.Range("J4").Formula = "=B4-SUM(C4:I4)"
$ws.Range("J5").Formula = "=B5-SUM(C5:I5)"
$ws.Range("J6").Formula = "=B6-SUM(C6:I6)"
$ws.Range("J7").Formula = "=B7-SUM(C7:I7)"
$ws.Range("J8").Formula = "=B8-SUM(C8:I8)"
$ws.Range("J9").Formula = "=B9-SUM(C9:I9)"
$ws.Range("J10").Formula = "=B10-SUM(C10:I10)"
$ws.Range("J11").Formula = "=B11-SUM(C11:I11)"
$ws.Range("J12").Formula = "=B12-SUM(C12:I12)"
$ws.Range("J13").Formula = "=B13-SUM(C13:I13)"
$ws.Range("J14").Formula = "=B14-SUM(C14:I14)"
$ws.Range("J15").Formula = "=B15-SUM(C15:I15)"
$ws.Range("J16").Formula = "=B16-SUM(C16:I16)"
$ws.Range("J17").Formula = "=B17-SUM(C17:I17)"
$ws.Range("J18").Formula = "=B18-SUM(C18:I18)"
$ws.Range("J19").Formula = "=B19-SUM(C19:I19)"

# --- 5. Highlight the latest "Estimado" data point (I21) like the weekly header cells ---
$ws.Range("I21").Interior.Color = 65535
$ws.Range("I21").NumberFormat = "m/d/yyyy"

# --- 6. Reposition / resize the burndown chart now that the sheet is narrower ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 660.8940440452756
$co.Top = 8.69740157480315
$co.Width = 1298.9045275590552
$co.Height = 559.9633070866141

# --- 7. View tweaks made while updating the sheet ---
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$excel.ActiveWindow.Zoom = 80
$ws.Range("AL43").Select() | Out-Null
